# Apply updated crypto price/volume data to Sheet1
# (values are forced to Text via a leading apostrophe so Excel does not
#  reinterpret dotted price strings like "27.951.83" as numbers, then the
#  cell style is reset to "Normal" so no extra quote-prefix formatting is kept)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2: D2: '27.946.65' -> '27.951.83'; E2: '  -0.51%  ' -> '  -0.44%  '
Set-TextValue "D2" "27.951.83"
Set-TextValue "E2" "  -0.44%  "

# Row 3: D3: '1.857.72' -> '1.857.49'; E3: '  -0.96%  ' -> '  -0.94%  '
Set-TextValue "D3" "1.857.49"
Set-TextValue "E3" "  -0.94%  "

# Row 4: E4: '  -0.05%  ' -> '  -0.06%  '
Set-TextValue "E4" "  -0.06%  "

# Row 5: D5: '311.96' -> '311.99'; E5: '  -0.53%  ' -> '  -0.44%  '
Set-TextValue "D5" "311.99"
Set-TextValue "E5" "  -0.44%  "

# Row 6: E6: '  -0.05%  ' -> '  -0.04%  '
Set-TextValue "E6" "  -0.04%  "

# Row 7: D7: '0.5138' -> '0.5137'; E7: '  +1.70%  ' -> '  +1.85%  '
Set-TextValue "D7" "0.5137"
Set-TextValue "E7" "  +1.85%  "

# Row 8: D8: '0.3825' -> '0.3821'; E8: '  -0.24%  ' -> '  -0.25%  '
Set-TextValue "D8" "0.3821"
Set-TextValue "E8" "  -0.25%  "

# Row 9: D9: '0.08233' -> '0.08227'; E9: '  -4.59%  ' -> '  -4.62%  '
Set-TextValue "D9" "0.08227"
Set-TextValue "E9" "  -4.62%  "

# Row 10: D10: '1.109' -> '1.108'
Set-TextValue "D10" "1.108"

# Row 11: D11: '41.43' -> '41.45'; E11: '  -0.13%  ' -> '  -0.07%  '
Set-TextValue "D11" "41.45"
Set-TextValue "E11" "  -0.07%  "

# Row 12: E12: '  -2.52%  ' -> '  -2.47%  '
Set-TextValue "E12" "  -2.47%  "

# Row 13: B13: 'Solana' -> 'WrappedEther'; C13: 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol' -> 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D13: '20.49' -> '1.866.11'; E13: '  -0.87%  ' -> '  -0.13%  '
Set-TextValue "B13" "WrappedEther"
Set-TextValue "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.866.11"
Set-TextValue "E13" "  -0.13%  "

# Row 14: B14: 'WrappedEther' -> 'Solana'; C14: 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' -> 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D14: '1.861.33' -> '20.47'; E14: '  -0.63%  ' -> '  -1.01%  '
Set-TextValue "B14" "Solana"
Set-TextValue "C14" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D14" "20.47"
Set-TextValue "E14" "  -1.01%  "

# Row 15: D15: '7.245' -> '7.247'; E15: '  +1.05%  ' -> '  +1.16%  '
Set-TextValue "D15" "7.247"
Set-TextValue "E15" "  +1.16%  "

# Row 16: E16: '  -0.13%  ' -> '  -0.11%  '
Set-TextValue "E16" "  -0.11%  "

# Row 17: D17: '0.00001096' -> '0.00001095'; E17: '  -0.49%  ' -> '  -0.60%  '
Set-TextValue "D17" "0.00001095"
Set-TextValue "E17" "  -0.60%  "

# Row 18: D18: '90.32' -> '90.30'; E18: '  -0.83%  ' -> '  -0.79%  '
Set-TextValue "D18" "90.30"
Set-TextValue "E18" "  -0.79%  "

# Row 19: E19: '  +0.04%  ' -> '  +0.13%  '
Set-TextValue "E19" "  +0.13%  "

# Row 20: D20: '17.65' -> '17.64'; E20: '  -2.68%  ' -> '  -2.60%  '
Set-TextValue "D20" "17.64"
Set-TextValue "E20" "  -2.60%  "

# Row 21: E21: '  -0.10%  ' -> '  -0.09%  '
Set-TextValue "E21" "  -0.09%  "

# Row 22: D22: '6.002' -> '6.004'; E22: '  -1.64%  ' -> '  -1.58%  '
Set-TextValue "D22" "6.004"
Set-TextValue "E22" "  -1.58%  "

# Row 23: D23: '27.983.12' -> '27.986.49'; E23: '  -0.53%  ' -> '  -0.47%  '
Set-TextValue "D23" "27.986.49"
Set-TextValue "E23" "  -0.47%  "

# Row 24: E24: '  -3.53%  ' -> '  -3.49%  '
Set-TextValue "E24" "  -3.49%  "

# Row 25: E25: '  -1.45%  ' -> '  -1.23%  '
Set-TextValue "E25" "  -1.23%  "

# Row 26: D26: '2.072.17' -> '2.071.91'; E26: '  -0.72%  ' -> '  -0.63%  '
Set-TextValue "D26" "2.071.91"
Set-TextValue "E26" "  -0.63%  "

# Row 27: D27: '2.501' -> '2.499'; E27: '  -2.80%  ' -> '  -2.48%  '
Set-TextValue "D27" "2.499"
Set-TextValue "E27" "  -2.48%  "

# Row 28: D28: '157.26' -> '157.29'; E28: '  +0.05%  ' -> '  +0.07%  '
Set-TextValue "D28" "157.29"
Set-TextValue "E28" "  +0.07%  "

# Row 29: E29: '  -1.44%  ' -> '  -1.45%  '
Set-TextValue "E29" "  -1.45%  "

# Row 30: D30: '124.48' -> '124.49'; E30: '  -1.21%  ' -> '  -1.15%  '
Set-TextValue "D30" "124.49"
Set-TextValue "E30" "  -1.15%  "

# Row 31: E31: '  +1.28%  ' -> '  +1.24%  '
Set-TextValue "E31" "  +1.24%  "

# Row 32: D32: '1.025' -> '1.024'; E32: '  -3.27%  ' -> '  -3.36%  '
Set-TextValue "D32" "1.024"
Set-TextValue "E32" "  -3.36%  "

# Row 33: D33: '5.862' -> '5.843'; E33: '  +4.86%  ' -> '  +4.54%  '
Set-TextValue "D33" "5.843"
Set-TextValue "E33" "  +4.54%  "

# Row 34: D34: '3.590' -> '3.592'; E34: '  -0.48%  ' -> '  -0.52%  '
Set-TextValue "D34" "3.592"
Set-TextValue "E34" "  -0.52%  "

# Row 35: D35: '9.395' -> '9.372'; E35: '  -3.03%  ' -> '  -3.11%  '
Set-TextValue "D35" "9.372"
Set-TextValue "E35" "  -3.11%  "

# Row 36: E36: '  -1.79%  ' -> '  -1.73%  '
Set-TextValue "E36" "  -1.73%  "

# Row 37: D37: '0.06497' -> '0.06490'; E37: '  -1.22%  ' -> '  -1.25%  '
Set-TextValue "D37" "0.06490"
Set-TextValue "E37" "  -1.25%  "

# Row 38: E38: '  +0.35%  ' -> '  +0.28%  '
Set-TextValue "E38" "  +0.28%  "

# Row 39: D39: '0.6534' -> '0.6537'; E39: '  +2.54%  ' -> '  +2.75%  '
Set-TextValue "D39" "0.6537"
Set-TextValue "E39" "  +2.75%  "

# Row 40: D40: '1.192' -> '1.191'; E40: '  -1.10%  ' -> '  -1.09%  '
Set-TextValue "D40" "1.191"
Set-TextValue "E40" "  -1.09%  "

# Row 42: E42: '  -3.14%  ' -> '  -3.26%  '
Set-TextValue "E42" "  -3.26%  "

# Row 43: D43: '11.14' -> '11.13'; E43: '  -3.48%  ' -> '  -3.86%  '
Set-TextValue "D43" "11.13"
Set-TextValue "E43" "  -3.86%  "

# Row 44: D44: '0.6109' -> '0.6103'; E44: '  +2.12%  ' -> '  +2.01%  '
Set-TextValue "D44" "0.6103"
Set-TextValue "E44" "  +2.01%  "

# Row 45: D45: '13.01' -> '12.93'; E45: '  -1.97%  ' -> '  -1.65%  '
Set-TextValue "D45" "12.93"
Set-TextValue "E45" "  -1.65%  "

# Row 46: D46: '3.674' -> '3.677'; E46: '  +0.01%  ' -> '  +0.15%  '
Set-TextValue "D46" "3.677"
Set-TextValue "E46" "  +0.15%  "

# Row 47: D47: '1.273' -> '1.272'; E47: '  -0.54%  ' -> '  -0.63%  '
Set-TextValue "D47" "1.272"
Set-TextValue "E47" "  -0.63%  "

# Row 48: D48: '2.006' -> '2.005'; E48: '  +0.80%  ' -> '  +0.81%  '
Set-TextValue "D48" "2.005"
Set-TextValue "E48" "  +0.81%  "

# Row 49: E49: '  -1.67%  ' -> '  -1.53%  '
Set-TextValue "E49" "  -1.53%  "

# Row 50: E50: '  -0.55%  ' -> '  -0.47%  '
Set-TextValue "E50" "  -0.47%  "

# Row 51: D51: '78.13' -> '78.08'; E51: '  -2.62%  ' -> '  -2.66%  '
Set-TextValue "D51" "78.08"
Set-TextValue "E51" "  -2.66%  "
